$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws "D2" "27.605.66"
Set-TextCell $ws "E2" "  +0.31%  "
Set-TextCell $ws "D3" "1.846.68"
Set-TextCell $ws "E3" "  +0.40%  "
Set-TextCell $ws "D4" "1.031"
Set-TextCell $ws "E4" "  +0.28%  "
Set-TextCell $ws "D5" "321.48"
Set-TextCell $ws "E5" "  +1.10%  "
Set-TextCell $ws "D6" "1.026"
Set-TextCell $ws "E6" "  +0.19%  "
Set-TextCell $ws "D7" "0.4373"
Set-TextCell $ws "E7" "  +0.00%  "
Set-TextCell $ws "E8" "  +1.34%  "
Set-TextCell $ws "D9" "0.07367"
Set-TextCell $ws "E9" "  -0.09%  "
Set-TextCell $ws "D10" "0.8801"
Set-TextCell $ws "E10" "  +0.52%  "
Set-TextCell $ws "D11" "21.49"
Set-TextCell $ws "E11" "  +0.09%  "
Set-TextCell $ws "D12" "1.854.67"
Set-TextCell $ws "D13" "5.494"
Set-TextCell $ws "E13" "  +0.05%  "
Set-TextCell $ws "D14" "6.696"
Set-TextCell $ws "E14" "  +0.32%  "
Set-TextCell $ws "D15" "0.07128"
Set-TextCell $ws "E15" "  +0.08%  "
Set-TextCell $ws "D16" "84.97"
Set-TextCell $ws "E16" "  +2.83%  "
Set-TextCell $ws "E17" "  +0.19%  "
Set-TextCell $ws "D18" "0.000009038"
Set-TextCell $ws "E18" "  +0.43%  "
Set-TextCell $ws "E19" "  +0.26%  "
Set-TextCell $ws "E20" "  +0.15%  "
Set-TextCell $ws "D21" "27.621.73"
Set-TextCell $ws "E21" "  +0.30%  "
Set-TextCell $ws "D22" "5.277"
Set-TextCell $ws "E22" "  +0.90%  "
Set-TextCell $ws "D23" "11.28"
Set-TextCell $ws "E23" "  +0.77%  "
Set-TextCell $ws "D24" "2.078.83"
Set-TextCell $ws "E24" "  +0.42%  "
Set-TextCell $ws "D25" "2.033"
Set-TextCell $ws "E25" "  +6.24%  "
Set-TextCell $ws "D26" "157.46"
Set-TextCell $ws "E26" "  +0.44%  "
Set-TextCell $ws "D27" "18.66"
Set-TextCell $ws "D28" "1.999"
Set-TextCell $ws "E28" "  +3.27%  "
Set-TextCell $ws "D29" "5.319"
Set-TextCell $ws "E29" "  +1.31%  "
Set-TextCell $ws "D30" "117.65"
Set-TextCell $ws "E30" "  +1.21%  "
Set-TextCell $ws "D31" "0.09017"
Set-TextCell $ws "E31" "  -0.49%  "
Set-TextCell $ws "B32" "ImmutableX"
Set-TextCell $ws "C32" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell $ws "D32" "0.7687"
Set-TextCell $ws "E32" "  +0.39%  "
Set-TextCell $ws "B33" "ARBITRUM"
Set-TextCell $ws "C33" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell $ws "D33" "1.205"
Set-TextCell $ws "E33" "  -0.36%  "
Set-TextCell $ws "D34" "2.987"
Set-TextCell $ws "E34" "  +4.39%  "
Set-TextCell $ws "E35" "  +0.95%  "
Set-TextCell $ws "D36" "1.027"
Set-TextCell $ws "E36" "  +0.17%  "
Set-TextCell $ws "D37" "1.141"
Set-TextCell $ws "E37" "  -0.41%  "
Set-TextCell $ws "D38" "0.01967"
Set-TextCell $ws "E38" "  -0.21%  "
Set-TextCell $ws "D39" "0.05259"
Set-TextCell $ws "E39" "  +0.05%  "
Set-TextCell $ws "D40" "2.839"
Set-TextCell $ws "E40" "  +1.33%  "
Set-TextCell $ws "D41" "0.5165"
Set-TextCell $ws "E41" "  -0.06%  "
Set-TextCell $ws "E42" "  -0.08%  "
Set-TextCell $ws "D43" "6.822"
Set-TextCell $ws "E43" "  +3.18%  "
Set-TextCell $ws "D44" "8.743"
Set-TextCell $ws "E44" "  +2.36%  "
Set-TextCell $ws "D45" "110.17"
Set-TextCell $ws "D46" "10.68"
Set-TextCell $ws "E46" "  +1.10%  "
Set-TextCell $ws "D47" "0.06605"
Set-TextCell $ws "E47" "  +4.25%  "
Set-TextCell $ws "D48" "1.028"
Set-TextCell $ws "E48" "  +0.37%  "
Set-TextCell $ws "D49" "1.697"
Set-TextCell $ws "E49" "  -0.34%  "
Set-TextCell $ws "D50" "0.4683"
Set-TextCell $ws "E50" "  +0.81%  "
Set-TextCell $ws "D51" "1.884"
Set-TextCell $ws "E51" "  -1.13%  "
